$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. _input sheet: insert a new row 4 ("数据粒度" / "按月查看") which
#    pushes the old row4 ("时间"/date) to row5 and old row5
#    ("高级经理"/amount) to row6.
# ---------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("_input")
$wsInput.Rows.Item(4).Insert()

$wsInput.Range("A4").Value = "数据粒度"
$wsInput.Range("B4").Value = "按月查看"
$wsInput.Rows.Item(4).RowHeight = 18

# Row 3 no longer needs the "thick bottom border" row hint now that the
# divider sits one row further down (new row 4's bottom edge takes over).
$wsInput.Rows.Item(3).AutoFit()
$wsInput.Rows.Item(3).RowHeight = 18

# ---------------------------------------------------------------
# 2. trend sheet: title formula now also pulls in the new
#    "数据粒度" value cell (_input!$B4) and calls the result a
#    report ("报表") instead of a trend chart ("趋势图").
# ---------------------------------------------------------------
$wsTrend = $wb.Worksheets.Item("trend")
$wsTrend.Range("B2").Formula = '=_input!$B2&_input!$B3&_input!$B4&"报表"'

# ---------------------------------------------------------------
# 3. Workbook-level defined name: shift the OFFSET anchor down one
#    row and bump the row-count correction to account for the
#    extra header row in _input.
# ---------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*trendChartDataSource*") {
        $n.RefersTo = "=OFFSET(_input!`$A5,0,0,COUNTA(_input!`$A:`$A)-4, COUNTA(_input!`$5:`$5))"
    }
}

# ---------------------------------------------------------------
# 4. Chart series: point at the shifted _input cells.
# ---------------------------------------------------------------
$co = $wsTrend.ChartObjects().Item(1)
$cht = $co.Chart
$ser = $cht.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(_input!`$A`$6,_input!`$B`$5:`$B`$5,_input!`$B`$6:`$B`$6,1)"
